$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set new column width for column 55 (BC) to match width="12" ---

# --- BC1: header date text "2024/11/02" (must stay literal text, not auto-converted to a date) ---
$ws.Range("ZZ1").Formula = '="2024/11/02"'
$ws.Range("A1").Copy()
$ws.Range("BC1").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("BC1").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

# --- BC2:BC53: numeric values with per-row fill style copied from a reference cell ---
$ws.Range("N2").Copy()
$ws.Range("BC2").PasteSpecial(-4122)
$ws.Range("BC2").Value = 131.1
$ws.Range("A2").Copy()
$ws.Range("BC3").PasteSpecial(-4122)
$ws.Range("BC3").Value = 195.2
$ws.Range("D2").Copy()
$ws.Range("BC4").PasteSpecial(-4122)
$ws.Range("BC4").Value = 115.1
$ws.Range("A2").Copy()
$ws.Range("BC5").PasteSpecial(-4122)
$ws.Range("BC5").Value = 186
$ws.Range("N2").Copy()
$ws.Range("BC6").PasteSpecial(-4122)
$ws.Range("BC6").Value = 135.8
$ws.Range("A2").Copy()
$ws.Range("BC7").PasteSpecial(-4122)
$ws.Range("BC7").Value = 178.5
$ws.Range("A2").Copy()
$ws.Range("BC8").PasteSpecial(-4122)
$ws.Range("BC8").Value = 171.8
$ws.Range("N2").Copy()
$ws.Range("BC9").PasteSpecial(-4122)
$ws.Range("BC9").Value = 134
$ws.Range("A2").Copy()
$ws.Range("BC10").PasteSpecial(-4122)
$ws.Range("BC10").Value = 196.5
$ws.Range("A2").Copy()
$ws.Range("BC11").PasteSpecial(-4122)
$ws.Range("BC11").Value = 168.6
$ws.Range("A2").Copy()
$ws.Range("BC12").PasteSpecial(-4122)
$ws.Range("BC12").Value = 183
$ws.Range("A2").Copy()
$ws.Range("BC13").PasteSpecial(-4122)
$ws.Range("BC13").Value = 144.4
$ws.Range("A2").Copy()
$ws.Range("BC14").PasteSpecial(-4122)
$ws.Range("BC14").Value = 149
$ws.Range("A2").Copy()
$ws.Range("BC15").PasteSpecial(-4122)
$ws.Range("BC15").Value = 149.4
$ws.Range("A2").Copy()
$ws.Range("BC16").PasteSpecial(-4122)
$ws.Range("BC16").Value = 155
$ws.Range("D2").Copy()
$ws.Range("BC17").PasteSpecial(-4122)
$ws.Range("BC17").Value = 117.4
$ws.Range("A2").Copy()
$ws.Range("BC18").PasteSpecial(-4122)
$ws.Range("BC18").Value = 172
$ws.Range("A2").Copy()
$ws.Range("BC19").PasteSpecial(-4122)
$ws.Range("BC19").Value = 168.6
$ws.Range("A2").Copy()
$ws.Range("BC20").PasteSpecial(-4122)
$ws.Range("BC20").Value = 146.5
$ws.Range("A2").Copy()
$ws.Range("BC21").PasteSpecial(-4122)
$ws.Range("BC21").Value = 200.1
$ws.Range("A2").Copy()
$ws.Range("BC22").PasteSpecial(-4122)
$ws.Range("BC22").Value = 205.9
$ws.Range("N2").Copy()
$ws.Range("BC23").PasteSpecial(-4122)
$ws.Range("BC23").Value = 135.3
$ws.Range("N2").Copy()
$ws.Range("BC24").PasteSpecial(-4122)
$ws.Range("BC24").Value = 138.9
$ws.Range("A2").Copy()
$ws.Range("BC25").PasteSpecial(-4122)
$ws.Range("BC25").Value = 160.6
$ws.Range("A2").Copy()
$ws.Range("BC26").PasteSpecial(-4122)
$ws.Range("BC26").Value = 150.7
$ws.Range("A2").Copy()
$ws.Range("BC27").PasteSpecial(-4122)
$ws.Range("BC27").Value = 150.3
$ws.Range("N2").Copy()
$ws.Range("BC28").PasteSpecial(-4122)
$ws.Range("BC28").Value = 131.7
$ws.Range("N2").Copy()
$ws.Range("BC29").PasteSpecial(-4122)
$ws.Range("BC29").Value = 139.9
$ws.Range("A2").Copy()
$ws.Range("BC30").PasteSpecial(-4122)
$ws.Range("BC30").Value = 235.4
$ws.Range("A2").Copy()
$ws.Range("BC31").PasteSpecial(-4122)
$ws.Range("BC31").Value = 222.6
$ws.Range("A2").Copy()
$ws.Range("BC32").PasteSpecial(-4122)
$ws.Range("BC32").Value = 141.2
$ws.Range("D2").Copy()
$ws.Range("BC33").PasteSpecial(-4122)
$ws.Range("BC33").Value = 110.1
$ws.Range("N2").Copy()
$ws.Range("BC34").PasteSpecial(-4122)
$ws.Range("BC34").Value = 134.4
$ws.Range("A2").Copy()
$ws.Range("BC35").PasteSpecial(-4122)
$ws.Range("BC35").Value = 143
$ws.Range("A2").Copy()
$ws.Range("BC36").PasteSpecial(-4122)
$ws.Range("BC36").Value = 171.1
$ws.Range("N2").Copy()
$ws.Range("BC37").PasteSpecial(-4122)
$ws.Range("BC37").Value = 132.8
$ws.Range("A2").Copy()
$ws.Range("BC38").PasteSpecial(-4122)
$ws.Range("BC38").Value = 182.7
$ws.Range("N2").Copy()
$ws.Range("BC39").PasteSpecial(-4122)
$ws.Range("BC39").Value = 133.2
$ws.Range("A2").Copy()
$ws.Range("BC40").PasteSpecial(-4122)
$ws.Range("BC40").Value = 205.2
$ws.Range("A2").Copy()
$ws.Range("BC41").PasteSpecial(-4122)
$ws.Range("BC41").Value = 152.8
$ws.Range("A2").Copy()
$ws.Range("BC42").PasteSpecial(-4122)
$ws.Range("BC42").Value = 158.8
$ws.Range("A2").Copy()
$ws.Range("BC43").PasteSpecial(-4122)
$ws.Range("BC43").Value = 226.7
$ws.Range("D2").Copy()
$ws.Range("BC44").PasteSpecial(-4122)
$ws.Range("BC44").Value = 113.1
$ws.Range("A2").Copy()
$ws.Range("BC45").PasteSpecial(-4122)
$ws.Range("BC45").Value = 148.2
$ws.Range("A2").Copy()
$ws.Range("BC46").PasteSpecial(-4122)
$ws.Range("BC46").Value = 152.5
$ws.Range("A2").Copy()
$ws.Range("BC47").PasteSpecial(-4122)
$ws.Range("BC47").Value = 147.4
$ws.Range("A2").Copy()
$ws.Range("BC48").PasteSpecial(-4122)
$ws.Range("BC48").Value = 156.3
$ws.Range("N2").Copy()
$ws.Range("BC49").PasteSpecial(-4122)
$ws.Range("BC49").Value = 139.3
$ws.Range("N2").Copy()
$ws.Range("BC50").PasteSpecial(-4122)
$ws.Range("BC50").Value = 130.8
$ws.Range("A2").Copy()
$ws.Range("BC51").PasteSpecial(-4122)
$ws.Range("BC51").Value = 209.4
$ws.Range("A2").Copy()
$ws.Range("BC52").PasteSpecial(-4122)
$ws.Range("BC52").Value = 175.4
$ws.Range("N2").Copy()
$ws.Range("BC53").PasteSpecial(-4122)
$ws.Range("BC53").Value = 126.6

# --- set width of new column (55) to match target col width=12 ---
$ws.Range("BC1").ColumnWidth = 11.17

$excel.CutCopyMode = 0
